$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 350
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -130
$ws.Range("N12").Value = -840

$ws.Range("H113").Value = 2871.9565
$ws.Range("I113").Value = 2866.3333
$ws.Range("J113").Value = 2882.5
$ws.Range("K113").Value = 2866.3333
$ws.Range("L113").Value = 2882.5
$ws.Range("M113").Value = 387.6667000000002
$ws.Range("N113").Value = -9390.5

$ws.Range("H137").Value = 1149.875
$ws.Range("I137").Value = 950
$ws.Range("K137").Value = 2850
$ws.Range("M137").Value = -300


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1730.6957
$ws.Range("I2").Value = 989.1053000000001
$ws.Range("K2").Value = 989.1053000000001
$ws.Range("M2").Value = -876.1053000000001

$ws.Range("H32").Value = 320056.06
$ws.Range("I32").Value = 2264.4321
$ws.Range("K32").Value = 2264.4321
$ws.Range("M32").Value = -1977.4321

$ws.Range("H38").Value = 3125
$ws.Range("I38").Value = 2000
$ws.Range("K38").Value = 2000
$ws.Range("M38").Value = -1533

$ws.Range("H110").Value = 1426.8572
$ws.Range("I110").Value = 1331.3334
$ws.Range("K110").Value = 1331.3334
$ws.Range("M110").Value = 713.6666

$ws.Range("H116").Value = 1730.6957
$ws.Range("I116").Value = 989.1053000000001
$ws.Range("K116").Value = 989.1053000000001
$ws.Range("M116").Value = 1304.8947

$ws.Range("H122").Value = 19232690
$ws.Range("I122").Value = 21741094
$ws.Range("K122").Value = 65223282
$ws.Range("M122").Value = -65220832

$ws.Range("H132").Value = 11114178
$ws.Range("I132").Value = 1796.3334
$ws.Range("J132").Value = 18522434
$ws.Range("K132").Value = 5389.0002
$ws.Range("L132").Value = 55567302
$ws.Range("M132").Value = -2859.0002
$ws.Range("N132").Value = -55572362


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1730.6957
$ws.Range("I3").Value = 989.1053000000001
$ws.Range("K3").Value = 989.1053000000001
$ws.Range("M3").Value = -875.1053000000001

$ws.Range("H22").Value = 440
$ws.Range("I22").Value = 425
$ws.Range("K22").Value = 425
$ws.Range("M22").Value = -252


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7249095
$ws.Range("I31").Value = 10639712
$ws.Range("J31").Value = 5505.3184
$ws.Range("K31").Value = 10639712
$ws.Range("L31").Value = 5505.3184
$ws.Range("M31").Value = -10639417
$ws.Range("N31").Value = -6095.3184

$ws.Range("H34").Value = 7249095
$ws.Range("I34").Value = 10639712
$ws.Range("J34").Value = 5505.3184
$ws.Range("K34").Value = 10639712
$ws.Range("L34").Value = 5505.3184
$ws.Range("M34").Value = -10639510
$ws.Range("N34").Value = -5909.3184

$ws.Range("H132").Value = 28171.297
$ws.Range("I132").Value = 651.6070999999999
$ws.Range("J132").Value = 113788.11
$ws.Range("K132").Value = 1954.8213
$ws.Range("L132").Value = 341364.33
$ws.Range("M132").Value = 575.1787000000002
$ws.Range("N132").Value = -346424.33

$ws.Range("H141").Value = 19806.25
$ws.Range("I141").Value = 10000
$ws.Range("J141").Value = 21207.143
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 21207.143
$ws.Range("M141").Value = -4820
$ws.Range("N141").Value = -31567.143


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 980
$ws.Range("I113").Value = 723.625
$ws.Range("J113").Value = 1027.6976
$ws.Range("K113").Value = 2170.875
$ws.Range("L113").Value = 3083.0928
$ws.Range("M113").Value = -0.875
$ws.Range("N113").Value = -7423.0928

$ws.Range("H121").Value = 38979.58
$ws.Range("I121").Value = 90
$ws.Range("J121").Value = 44052.13
$ws.Range("K121").Value = 270
$ws.Range("L121").Value = 132156.39
$ws.Range("M121").Value = 1040
$ws.Range("N121").Value = -134776.39

$ws.Range("H131").Value = 7353964.5
$ws.Range("J131").Value = 9616542
$ws.Range("L131").Value = 28849626
$ws.Range("N131").Value = -28859706

$ws.Range("H132").Value = 2126.7058
$ws.Range("I132").Value = 1638.8889
$ws.Range("J132").Value = 2675.5
$ws.Range("K132").Value = 14750.0001
$ws.Range("L132").Value = 24079.5
$ws.Range("M132").Value = -12220.0001
$ws.Range("N132").Value = -29139.5

$ws.Range("H137").Value = 2735.0527
$ws.Range("I137").Value = 1400
$ws.Range("J137").Value = 2892.1177
$ws.Range("K137").Value = 4200
$ws.Range("L137").Value = 8676.3531
$ws.Range("M137").Value = 900
$ws.Range("N137").Value = -18876.3531


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 20083.334
$ws.Range("I46").Value = 9125
$ws.Range("J46").Value = 42000
$ws.Range("K46").Value = 9125
$ws.Range("L46").Value = 42000
$ws.Range("M46").Value = -8969
$ws.Range("N46").Value = -42312


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 1377.9
$ws.Range("I31").Value = 476
$ws.Range("J31").Value = 2279.8
$ws.Range("K31").Value = 476
$ws.Range("L31").Value = 2279.8
$ws.Range("M31").Value = -228
$ws.Range("N31").Value = -2775.8

$ws.Range("H61").Value = 1950.5714
$ws.Range("I61").Value = 1556.4445
$ws.Range("J61").Value = 2660
$ws.Range("K61").Value = 1556.4445
$ws.Range("L61").Value = 2660
$ws.Range("M61").Value = -1354.4445
$ws.Range("N61").Value = -3064

$ws.Range("H113").Value = 1950.5714
$ws.Range("I113").Value = 1556.4445
$ws.Range("J113").Value = 2660
$ws.Range("K113").Value = 1556.4445
$ws.Range("L113").Value = 2660
$ws.Range("M113").Value = 613.5554999999999
$ws.Range("N113").Value = -7000

$ws.Range("H122").Value = 1872.3864
$ws.Range("I122").Value = 1792
$ws.Range("J122").Value = 2185
$ws.Range("K122").Value = 5376
$ws.Range("L122").Value = 6555
$ws.Range("M122").Value = -2926
$ws.Range("N122").Value = -11455

$ws.Range("H136").Value = 2792.976
$ws.Range("I136").Value = 2016.6129
$ws.Range("J136").Value = 4980.909
$ws.Range("K136").Value = 6049.8387
$ws.Range("L136").Value = 14942.727
$ws.Range("M136").Value = -3499.8387
$ws.Range("N136").Value = -20042.727


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 41668532
$ws.Range("I126").Value = 71430190
$ws.Range("J126").Value = 2205.9
$ws.Range("K126").Value = 214290570
$ws.Range("L126").Value = 6617.700000000001
$ws.Range("M126").Value = -214288100
$ws.Range("N126").Value = -11557.7

